# Apply the edits described by the commit:
#   - D4: 3 -> 2
#   - F4: -2 -> -3
#   - H4: 36 -> 46
#   - active selection moves from E5 to D5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 2
$ws.Range("F4").Value = -3
$ws.Range("H4").Value = 46

$ws.Range("D5").Select()
